# Sample Project / Main.xlsx — "Rules" sheet, cell B11 (4th rule row, "Rule" column)
# changes from the shared string "R40" to the literal text "1".
#
# A plain `$cell.Value = "1"` would let Excel's smart-entry parsing treat
# "1" as a number (since it is unformatted/General), which is not what the
# source workbook shows (it stays a text cell, same cell style). Routing the
# literal text through a formula + copy/paste-values keeps the result as a
# genuine text value (not a number, not a quote-prefixed/re-styled cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
